$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "id"
$ws.Range("B6").Value = "nombre"
$ws.Range("C6").Value = "edad"

$ws.Range("A7").Value = 0
$ws.Range("B7").Value = "Jose"
$ws.Range("C7").Value = 29

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Sofia"
$ws.Range("C8").Value = 24

$ws.Range("A9").Value = 2
$ws.Range("B9").Value = "Carlos"
$ws.Range("C9").Value = 36
